$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.802.65'
$ws.Range('E2').Value = '  -1.58%  '
$ws.Range('D3').Value = '3.382.89'
$ws.Range('E3').Value = '  -2.04%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '571.50'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.66%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '141.61'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.69%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = '3.384.01'
$ws.Range('E8').Value = '  -1.95%  '
$ws.Range('E9').Value = '  -0.08%  '
$ws.Range('E10').Value = '  -3.89%  '
$ws.Range('E11').Value = '  -0.67%  '
$ws.Range('E12').Value = '  +0.51%  '
$ws.Range('D13').Value = '3.961.37'
$ws.Range('E13').Value = '  -2.11%  '
$ws.Range('B14').Value = 'Avalanche'
$ws.Range('C14').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.10'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.04%  '
$ws.Range('B15').Value = 'TRON'
$ws.Range('C15').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.124'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.98%  '
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000171'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.63%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.386.26'
$ws.Range('E17').Value = '  -2.71%  '
$ws.Range('D18').Value = '60.894.56'
$ws.Range('E18').Value = '  -1.56%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.27'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.82%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.15'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.53%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.00'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.70%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '388.87'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.79%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.560'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.47%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '73.48'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.12%  '
$ws.Range('E25').Value = '  -0.18%  '
$ws.Range('E26').Value = '  -3.72%  '
$ws.Range('D27').Value = '3.523.89'
$ws.Range('E27').Value = '  -2.08%  '
$ws.Range('E28').Value = '  -0.68%  '
$ws.Range('E29').Value = '  -0.06%  '
$ws.Range('E30').Value = '  -5.05%  '
$ws.Range('E31').Value = '  -3.65%  '
$ws.Range('E32').Value = '  -2.02%  '
$ws.Range('E33').Value = '  -0.42%  '
$ws.Range('E34').Value = '  -0.05%  '
$ws.Range('E35').Value = '  -0.87%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.94'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.01%  '
$ws.Range('D37').Value = '3.413.12'
$ws.Range('E37').Value = '  -1.83%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '167.13'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.67%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.04'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.74%  '
$ws.Range('E40').Value = '  -3.25%  '
$ws.Range('E41').Value = '  -1.53%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '26.98'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.25%  '
$ws.Range('E43').Value = '  -1.78%  '
$ws.Range('E44').Value = '  -0.09%  '
$ws.Range('E45').Value = '  -0.72%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '41.74'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.48%  '
$ws.Range('E47').Value = '  -2.12%  '
$ws.Range('D48').Value = '2.536.96'
$ws.Range('E48').Value = '  -2.23%  '
$ws.Range('E49').Value = '  -4.17%  '
$ws.Range('E50').Value = '  -2.11%  '
$ws.Range('E51').Value = '  -1.43%  '
